$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new worksheet "corr0" right after "defaultvalues" (i.e.
#    right before "corr1"), holding the correlation matrix for the new
#    PARAM5 / PARAM6 correlated pair.
# ---------------------------------------------------------------------
$defaultvalues = $wb.Worksheets.Item("defaultvalues")
$corr0 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $defaultvalues)
$corr0.Name = "corr0"

$corr0.Range("B1").Value = "PARAM5"
$corr0.Range("C1").Value = "PARAM6"

$corr0.Range("A2").Value = "PARAM5"
$corr0.Range("B2").Value = 1

$corr0.Range("A3").Value = "PARAM6"
$corr0.Range("B3").Value = 0.8
$corr0.Range("C3").Value = 1

# ---------------------------------------------------------------------
# 2. Update the "designinput" sheet:
#    - sens6 (row 9) now samples 500 realizations instead of 10
#    - PARAM5 (row 9) and PARAM6 (row 10) are now correlated, so their
#      corr_sheet column (O) points at the new "corr0" sheet
# ---------------------------------------------------------------------
$designinput = $wb.Worksheets.Item("designinput")
$designinput.Range("B9").Value = 500
$designinput.Range("O9").Value = "corr0"
$designinput.Range("O10").Value = "corr0"

# ---------------------------------------------------------------------
# 3. Update view/selection state to match the authored workbook:
#    "corr0" becomes the active/selected sheet, with C8 selected; the
#    previously-active "designinput" sheet loses tabSelected and keeps
#    a new selection of B10.
# ---------------------------------------------------------------------
$null = $designinput.Range("B10").Select()

$corr0.Activate()
$null = $corr0.Range("C8").Select()
